# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# --- ALC row 29 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 923.0769   # H29: 928.5714 -> 923.0769
$ws.Cells.Item(29, 10).Value = 991.6667   # J29: 992.3077 -> 991.6667
$ws.Cells.Item(29, 12).Value = 2975.0001   # L29: 2976.9231 -> 2975.0001
$ws.Cells.Item(29, 14).Value = -3537.0001   # N29: -3538.9231 -> -3537.0001

# --- ALC row 74 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 7996.6665   # H74: 5957.2 -> 7996.6665
$ws.Cells.Item(74, 9).Value = 7996.6665   # I74: 5957.2 -> 7996.6665
$ws.Cells.Item(74, 11).Value = 7996.6665   # K74: 5957.2 -> 7996.6665
$ws.Cells.Item(74, 13).Value = -7060.6665   # M74: -5021.2 -> -7060.6665

# --- ALC row 77 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 7996.6665   # H77: 5957.2 -> 7996.6665
$ws.Cells.Item(77, 9).Value = 7996.6665   # I77: 5957.2 -> 7996.6665
$ws.Cells.Item(77, 11).Value = 39983.3325   # K77: 29786 -> 39983.3325
$ws.Cells.Item(77, 13).Value = -35303.3325   # M77: -25106 -> -35303.3325

# --- ALC row 88 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 757.6667   # H88: 749 -> 757.6667
$ws.Cells.Item(88, 10).Value = 757.6667   # J88: 749 -> 757.6667
$ws.Cells.Item(88, 12).Value = 757.6667   # L88: 749 -> 757.6667
$ws.Cells.Item(88, 14).Value = -1569.6667   # N88: -1561 -> -1569.6667

# --- ALC row 91 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(91, 8).Value = 757.6667   # H91: 749 -> 757.6667
$ws.Cells.Item(91, 10).Value = 757.6667   # J91: 749 -> 757.6667
$ws.Cells.Item(91, 12).Value = 757.6667   # L91: 749 -> 757.6667
$ws.Cells.Item(91, 14).Value = -3565.6667   # N91: -3557 -> -3565.6667

# --- ALC row 111 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 2164.125   # H111: 2023.1111 -> 2164.125
$ws.Cells.Item(111, 9).Value = 963.5   # I111: 953.7143 -> 963.5
$ws.Cells.Item(111, 11).Value = 2890.5   # K111: 2861.1429 -> 2890.5
$ws.Cells.Item(111, 13).Value = 176.5   # M111: 205.8571000000002 -> 176.5

# --- ALC row 132 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1237.2778   # H132: 1409.6428 -> 1237.2778
$ws.Cells.Item(132, 9).Value = 1237.2778   # I132: 1409.6428 -> 1237.2778
$ws.Cells.Item(132, 11).Value = 3711.8334   # K132: 4228.928400000001 -> 3711.8334
$ws.Cells.Item(132, 13).Value = -1181.8334   # M132: -1698.928400000001 -> -1181.8334

# --- ALC row 141 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 984   # H141: 1220.091 -> 984
$ws.Cells.Item(141, 9).Value = 984   # I141: 1243.6 -> 984
$ws.Cells.Item(141, 10).Value = 0   # J141: 985 -> 0
$ws.Cells.Item(141, 11).Value = 2952   # K141: 3730.8 -> 2952
$ws.Cells.Item(141, 12).Value = 0   # L141: 2955 -> 0
$ws.Cells.Item(141, 13).Value = 2228   # M141: 1449.2 -> 2228
$ws.Cells.Item(141, 14).ClearContents()   # N141: was -13315

# --- ARM row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6218.7295   # H32: 6560.485 -> 6218.7295
$ws.Cells.Item(32, 9).Value = 6218.7295   # I32: 6560.485 -> 6218.7295
$ws.Cells.Item(32, 11).Value = 6218.7295   # K32: 6560.485 -> 6218.7295
$ws.Cells.Item(32, 13).Value = -5931.7295   # M32: -6273.485 -> -5931.7295

# --- BSM row 7 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 13349   # H7: 13333.667 -> 13349
$ws.Cells.Item(7, 9).Value = 50   # I7: 0 -> 50
$ws.Cells.Item(7, 10).Value = 19998.5   # J7: 13333.667 -> 19998.5
$ws.Cells.Item(7, 11).Value = 50   # K7: 0 -> 50
$ws.Cells.Item(7, 12).Value = 19998.5   # L7: 13333.667 -> 19998.5
$ws.Cells.Item(7, 13).Value = 63   # M7: (new) -> 63
$ws.Cells.Item(7, 14).Value = -20224.5   # N7: -13559.667 -> -20224.5

# --- BSM row 20 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 6716.2856   # H20: 7684.5 -> 6716.2856
$ws.Cells.Item(20, 9).Value = 1403   # I20: 1527 -> 1403
$ws.Cells.Item(20, 11).Value = 1403   # K20: 1527 -> 1403
$ws.Cells.Item(20, 13).Value = -1156   # M20: -1280 -> -1156

# --- BSM row 88 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(88, 8).Value = 29939   # H88: 0 -> 29939
$ws.Cells.Item(88, 10).Value = 29939   # J88: 0 -> 29939
$ws.Cells.Item(88, 12).Value = 29939   # L88: 0 -> 29939
$ws.Cells.Item(88, 14).Value = -30751   # N88: (new) -> -30751

# --- BSM row 91 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(91, 8).Value = 29939   # H91: 0 -> 29939
$ws.Cells.Item(91, 10).Value = 29939   # J91: 0 -> 29939
$ws.Cells.Item(91, 12).Value = 29939   # L91: 0 -> 29939
$ws.Cells.Item(91, 14).Value = -32747   # N91: (new) -> -32747

# --- BSM row 94 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 4273.7144   # H94: 4026.375 -> 4273.7144
$ws.Cells.Item(94, 10).Value = 4574.75   # J94: 4118.8 -> 4574.75
$ws.Cells.Item(94, 12).Value = 4574.75   # L94: 4118.8 -> 4574.75
$ws.Cells.Item(94, 14).Value = -5476.75   # N94: -5020.8 -> -5476.75

# --- CRP row 3 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 401500   # H3: 287999.72 -> 401500
$ws.Cells.Item(3, 9).Value = 1000250   # I3: 501249.75 -> 1000250
$ws.Cells.Item(3, 10).Value = 2333.3333   # J3: 3666.3333 -> 2333.3333
$ws.Cells.Item(3, 11).Value = 1000250   # K3: 501249.75 -> 1000250
$ws.Cells.Item(3, 12).Value = 2333.3333   # L3: 3666.3333 -> 2333.3333
$ws.Cells.Item(3, 13).Value = -1000137   # M3: -501136.75 -> -1000137
$ws.Cells.Item(3, 14).Value = -2559.3333   # N3: -3892.3333 -> -2559.3333

# --- CRP row 132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3284   # H132: 3599.2 -> 3284
$ws.Cells.Item(132, 9).Value = 2623.25   # I132: 2666 -> 2623.25
$ws.Cells.Item(132, 10).Value = 4165   # J132: 4999 -> 4165
$ws.Cells.Item(132, 11).Value = 7869.75   # K132: 7998 -> 7869.75
$ws.Cells.Item(132, 12).Value = 12495   # L132: 14997 -> 12495
$ws.Cells.Item(132, 13).Value = -5339.75   # M132: -5468 -> -5339.75
$ws.Cells.Item(132, 14).Value = -17555   # N132: -20057 -> -17555

# --- CUL row 52 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(52, 8).Value = 0   # H52: 833 -> 0
$ws.Cells.Item(52, 10).Value = 0   # J52: 833 -> 0
$ws.Cells.Item(52, 12).Value = 0   # L52: 2499 -> 0
$ws.Cells.Item(52, 14).ClearContents()   # N52: was -3031

# --- CUL row 107 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 690   # H107: 751.5 -> 690
$ws.Cells.Item(107, 9).Value = 501.5   # I107: 503 -> 501.5
$ws.Cells.Item(107, 10).Value = 784.25   # J107: 801.2 -> 784.25
$ws.Cells.Item(107, 11).Value = 1504.5   # K107: 1509 -> 1504.5
$ws.Cells.Item(107, 12).Value = 2352.75   # L107: 2403.6 -> 2352.75
$ws.Cells.Item(107, 13).Value = 415.5   # M107: 411 -> 415.5
$ws.Cells.Item(107, 14).Value = -6192.75   # N107: -6243.6 -> -6192.75

# --- GSM row 11 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 6205999   # H11: 10000000 -> 6205999
$ws.Cells.Item(11, 9).Value = 7507500   # I11: 10000000 -> 7507500
$ws.Cells.Item(11, 10).Value = 999995   # J11: 0 -> 999995
$ws.Cells.Item(11, 11).Value = 7507500   # K11: 10000000 -> 7507500
$ws.Cells.Item(11, 12).Value = 999995   # L11: 0 -> 999995
$ws.Cells.Item(11, 13).Value = -7507361   # M11: -9999861 -> -7507361
$ws.Cells.Item(11, 14).Value = -1000273   # N11: (new) -> -1000273

# --- GSM row 102 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3118.5   # H102: 6699.5 -> 3118.5
$ws.Cells.Item(102, 9).Value = 3118.5   # I102: 6699.5 -> 3118.5
$ws.Cells.Item(102, 11).Value = 3118.5   # K102: 6699.5 -> 3118.5
$ws.Cells.Item(102, 13).Value = -1496.5   # M102: -5077.5 -> -1496.5

# --- GSM row 122 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2704.0833   # H122: 3352.7273 -> 2704.0833
$ws.Cells.Item(122, 9).Value = 3138.889   # I122: 3866 -> 3138.889
$ws.Cells.Item(122, 10).Value = 1399.6666   # J122: 2454.5 -> 1399.6666
$ws.Cells.Item(122, 11).Value = 9416.667000000001   # K122: 11598 -> 9416.667000000001
$ws.Cells.Item(122, 12).Value = 4198.9998   # L122: 7363.5 -> 4198.9998
$ws.Cells.Item(122, 13).Value = -6966.667000000001   # M122: -9148 -> -6966.667000000001
$ws.Cells.Item(122, 14).Value = -9098.9998   # N122: -12263.5 -> -9098.9998

# --- GSM row 132 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2153.9375   # H132: 2348.4285 -> 2153.9375
$ws.Cells.Item(132, 9).Value = 1497.5385   # I132: 1625.7273 -> 1497.5385
$ws.Cells.Item(132, 11).Value = 4492.6155   # K132: 4877.1819 -> 4492.6155
$ws.Cells.Item(132, 13).Value = -1962.6155   # M132: -2347.1819 -> -1962.6155

# --- LTW row 16 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 220   # H16: 259.5 -> 220
$ws.Cells.Item(16, 9).Value = 0   # I16: 299 -> 0
$ws.Cells.Item(16, 11).Value = 0   # K16: 299 -> 0
$ws.Cells.Item(16, 13).ClearContents()   # M16: was -129

# --- LTW row 22 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2597.7144   # H22: 2130.4443 -> 2597.7144
$ws.Cells.Item(22, 9).Value = 3110.818   # I22: 2292.375 -> 3110.818
$ws.Cells.Item(22, 10).Value = 716.3333   # J22: 835 -> 716.3333
$ws.Cells.Item(22, 11).Value = 3110.818   # K22: 2292.375 -> 3110.818
$ws.Cells.Item(22, 12).Value = 716.3333   # L22: 835 -> 716.3333
$ws.Cells.Item(22, 13).Value = -2815.818   # M22: -1997.375 -> -2815.818
$ws.Cells.Item(22, 14).Value = -1306.3333   # N22: -1425 -> -1306.3333

# --- LTW row 27 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 2597.7144   # H27: 2130.4443 -> 2597.7144
$ws.Cells.Item(27, 9).Value = 3110.818   # I27: 2292.375 -> 3110.818
$ws.Cells.Item(27, 10).Value = 716.3333   # J27: 835 -> 716.3333
$ws.Cells.Item(27, 11).Value = 3110.818   # K27: 2292.375 -> 3110.818
$ws.Cells.Item(27, 12).Value = 716.3333   # L27: 835 -> 716.3333
$ws.Cells.Item(27, 13).Value = -3003.818   # M27: -2185.375 -> -3003.818
$ws.Cells.Item(27, 14).Value = -930.3333   # N27: -1049 -> -930.3333

# --- LTW row 55 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1235.0476   # H55: 1184.4546 -> 1235.0476
$ws.Cells.Item(55, 9).Value = 2385   # I55: 2673 -> 2385
$ws.Cells.Item(55, 10).Value = 372.58334   # J55: 333.85715 -> 372.58334
$ws.Cells.Item(55, 11).Value = 2385   # K55: 2673 -> 2385
$ws.Cells.Item(55, 12).Value = 372.58334   # L55: 333.85715 -> 372.58334
$ws.Cells.Item(55, 13).Value = -2212   # M55: -2500 -> -2212
$ws.Cells.Item(55, 14).Value = -718.58334   # N55: -679.85715 -> -718.58334

# --- LTW row 82 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 3095.7144   # H82: 2618.889 -> 3095.7144
$ws.Cells.Item(82, 9).Value = 3095.7144   # I82: 2618.889 -> 3095.7144
$ws.Cells.Item(82, 11).Value = 3095.7144   # K82: 2618.889 -> 3095.7144
$ws.Cells.Item(82, 13).Value = -2734.7144   # M82: -2257.889 -> -2734.7144

# --- LTW row 85 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 3095.7144   # H85: 2618.889 -> 3095.7144
$ws.Cells.Item(85, 9).Value = 3095.7144   # I85: 2618.889 -> 3095.7144
$ws.Cells.Item(85, 11).Value = 3095.7144   # K85: 2618.889 -> 3095.7144
$ws.Cells.Item(85, 13).Value = -1847.7144   # M85: -1370.889 -> -1847.7144

# --- LTW row 93 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 3571   # H93: 3459.1667 -> 3571
$ws.Cells.Item(93, 10).Value = 3751   # J93: 3467.3333 -> 3751
$ws.Cells.Item(93, 12).Value = 3751   # L93: 3467.3333 -> 3751
$ws.Cells.Item(93, 14).Value = -6247   # N93: -5963.3333 -> -6247

# --- LTW row 132 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2748.037   # H132: 2981.2917 -> 2748.037
$ws.Cells.Item(132, 9).Value = 2220.05   # I132: 2453.1765 -> 2220.05
$ws.Cells.Item(132, 10).Value = 4256.5713   # J132: 4263.857 -> 4256.5713
$ws.Cells.Item(132, 11).Value = 6660.150000000001   # K132: 7359.529500000001 -> 6660.150000000001
$ws.Cells.Item(132, 12).Value = 12769.7139   # L132: 12791.571 -> 12769.7139
$ws.Cells.Item(132, 13).Value = -4130.150000000001   # M132: -4829.529500000001 -> -4130.150000000001
$ws.Cells.Item(132, 14).Value = -17829.7139   # N132: -17851.571 -> -17829.7139

# --- LTW row 136 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 4399.2   # H136: 4499.5 -> 4399.2
$ws.Cells.Item(136, 9).Value = 4399.2   # I136: 4499.5 -> 4399.2
$ws.Cells.Item(136, 11).Value = 13197.6   # K136: 13498.5 -> 13197.6
$ws.Cells.Item(136, 13).Value = -10647.6   # M136: -10948.5 -> -10647.6

# --- WVR row 14 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 3000   # H14: 3666.6667 -> 3000
$ws.Cells.Item(14, 9).Value = 3000   # I14: 3666.6667 -> 3000
$ws.Cells.Item(14, 11).Value = 3000   # K14: 3666.6667 -> 3000
$ws.Cells.Item(14, 13).Value = -2832   # M14: -3498.6667 -> -2832

# --- WVR row 97 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(97, 8).Value = 0   # H97: 21000 -> 0
$ws.Cells.Item(97, 10).Value = 0   # J97: 21000 -> 0
$ws.Cells.Item(97, 12).Value = 0   # L97: 21000 -> 0
$ws.Cells.Item(97, 14).ClearContents()   # N97: was -22982

# --- WVR row 122 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3768.2856   # H122: 2724.7 -> 3768.2856
$ws.Cells.Item(122, 9).Value = 4884   # I122: 2805.5293 -> 4884
$ws.Cells.Item(122, 10).Value = 1760   # J122: 2266.6667 -> 1760
$ws.Cells.Item(122, 11).Value = 14652   # K122: 8416.5879 -> 14652
$ws.Cells.Item(122, 12).Value = 5280   # L122: 6800.000100000001 -> 5280
$ws.Cells.Item(122, 13).Value = -12202   # M122: -5966.5879 -> -12202
$ws.Cells.Item(122, 14).Value = -10180   # N122: -11700.0001 -> -10180

# --- WVR row 126 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1482.1428   # H126: 1184.3 -> 1482.1428
$ws.Cells.Item(126, 9).Value = 1443.1818   # I126: 1030.4445 -> 1443.1818
$ws.Cells.Item(126, 10).Value = 1525   # J126: 1415.0834 -> 1525
$ws.Cells.Item(126, 11).Value = 4329.5454   # K126: 3091.3335 -> 4329.5454
$ws.Cells.Item(126, 12).Value = 4575   # L126: 4245.2502 -> 4575
$ws.Cells.Item(126, 13).Value = -1859.5454   # M126: -621.3335000000002 -> -1859.5454
$ws.Cells.Item(126, 14).Value = -9515   # N126: -9185.2502 -> -9515

# --- WVR row 132 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1084.8948   # H132: 1322.8 -> 1084.8948
$ws.Cells.Item(132, 9).Value = 628   # I132: 737 -> 628
$ws.Cells.Item(132, 10).Value = 2798.25   # J132: 3666 -> 2798.25
$ws.Cells.Item(132, 11).Value = 1884   # K132: 2211 -> 1884
$ws.Cells.Item(132, 12).Value = 8394.75   # L132: 10998 -> 8394.75
$ws.Cells.Item(132, 13).Value = 646   # M132: 319 -> 646
$ws.Cells.Item(132, 14).Value = -13454.75   # N132: -16058 -> -13454.75

# --- WVR row 136 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 2096.0344   # H136: 2159.25 -> 2096.0344
$ws.Cells.Item(136, 9).Value = 1895.7693   # I136: 1958.56 -> 1895.7693
$ws.Cells.Item(136, 11).Value = 5687.3079   # K136: 5875.68 -> 5687.3079
$ws.Cells.Item(136, 13).Value = -3137.3079   # M136: -3325.68 -> -3137.3079
